# This script applies the edit described by the diff:
# For a set of rows in column B (the "is_churn" label column), flip the
# stored value (0 <-> 1). The full list of (row, newValue) pairs was
# derived directly from the unified diff of the workbook's XML.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$changes = @(
    ,(2,0)
    ,(4,1)
    ,(57,1)
    ,(59,1)
    ,(88,0)
    ,(91,1)
    ,(93,1)
    ,(98,1)
    ,(107,0)
    ,(111,0)
    ,(121,1)
    ,(212,1)
    ,(220,0)
    ,(222,1)
    ,(279,0)
    ,(289,1)
    ,(302,1)
    ,(306,1)
    ,(318,0)
    ,(326,0)
    ,(345,0)
    ,(347,1)
    ,(353,0)
    ,(359,1)
    ,(361,0)
    ,(381,0)
    ,(388,1)
    ,(397,1)
    ,(427,1)
    ,(467,1)
    ,(496,1)
    ,(505,1)
    ,(511,1)
    ,(550,1)
    ,(551,1)
    ,(564,1)
    ,(575,1)
    ,(596,1)
    ,(620,1)
    ,(642,1)
    ,(645,1)
    ,(651,1)
    ,(668,0)
    ,(669,1)
    ,(726,0)
    ,(745,1)
    ,(779,1)
    ,(789,0)
    ,(791,1)
    ,(798,0)
    ,(816,1)
    ,(902,1)
    ,(914,1)
    ,(934,1)
    ,(937,1)
    ,(956,1)
    ,(963,0)
    ,(982,1)
    ,(986,1)
    ,(988,1)
    ,(994,1)
    ,(1009,1)
    ,(1026,1)
    ,(1038,1)
    ,(1041,1)
    ,(1079,0)
    ,(1091,1)
    ,(1097,1)
    ,(1113,0)
    ,(1126,1)
    ,(1147,1)
    ,(1159,1)
    ,(1183,1)
    ,(1191,0)
    ,(1235,0)
    ,(1240,0)
    ,(1264,1)
    ,(1275,1)
    ,(1285,1)
    ,(1349,1)
    ,(1376,1)
    ,(1411,1)
    ,(1432,0)
    ,(1455,1)
    ,(1503,1)
    ,(1513,0)
    ,(1518,0)
    ,(1543,0)
    ,(1557,0)
    ,(1558,1)
    ,(1574,1)
    ,(1605,1)
    ,(1626,1)
    ,(1627,1)
    ,(1641,1)
    ,(1643,1)
    ,(1646,0)
    ,(1653,0)
    ,(1659,1)
    ,(1669,1)
    ,(1715,1)
    ,(1740,1)
    ,(1745,0)
    ,(1746,1)
    ,(1752,1)
    ,(1761,0)
    ,(1775,0)
    ,(1807,1)
    ,(1870,1)
    ,(1875,1)
    ,(1878,0)
    ,(1881,0)
    ,(1890,0)
    ,(1895,1)
    ,(1927,0)
    ,(1932,1)
    ,(1935,1)
    ,(1937,1)
    ,(1942,1)
    ,(1964,0)
    ,(1970,1)
    ,(1975,1)
    ,(1981,1)
    ,(1997,0)
    ,(2009,0)
    ,(2045,0)
    ,(2063,0)
    ,(2084,0)
    ,(2090,1)
    ,(2097,1)
    ,(2099,0)
    ,(2108,1)
    ,(2111,1)
    ,(2112,1)
    ,(2122,1)
    ,(2130,1)
    ,(2160,1)
    ,(2204,1)
    ,(2211,0)
    ,(2233,1)
    ,(2246,1)
    ,(2248,0)
    ,(2250,0)
    ,(2255,0)
    ,(2262,0)
    ,(2270,1)
    ,(2283,1)
    ,(2288,1)
    ,(2290,1)
    ,(2294,1)
    ,(2310,1)
    ,(2329,1)
    ,(2330,1)
    ,(2345,1)
    ,(2379,0)
    ,(2390,1)
    ,(2399,0)
    ,(2435,1)
    ,(2458,0)
    ,(2489,1)
    ,(2507,1)
    ,(2532,1)
    ,(2571,0)
    ,(2575,1)
    ,(2602,1)
    ,(2604,1)
    ,(2609,0)
    ,(2628,0)
    ,(2636,0)
    ,(2639,1)
    ,(2649,1)
    ,(2666,1)
    ,(2685,1)
    ,(2701,0)
    ,(2732,0)
    ,(2738,0)
    ,(2744,1)
    ,(2779,0)
    ,(2790,1)
    ,(2819,0)
    ,(2832,1)
    ,(2834,0)
    ,(2858,1)
    ,(2869,1)
    ,(2876,1)
    ,(2917,1)
    ,(2918,0)
    ,(2936,1)
    ,(2950,1)
    ,(2952,1)
    ,(2976,1)
    ,(2982,1)
    ,(2987,1)
    ,(3051,1)
    ,(3059,0)
    ,(3066,1)
    ,(3067,0)
    ,(3088,1)
    ,(3090,1)
    ,(3137,1)
    ,(3164,1)
    ,(3171,1)
    ,(3200,1)
    ,(3236,1)
    ,(3270,0)
    ,(3293,1)
    ,(3298,1)
    ,(3316,1)
    ,(3333,1)
    ,(3367,1)
    ,(3387,0)
    ,(3420,0)
    ,(3446,0)
    ,(3488,0)
    ,(3504,1)
    ,(3508,0)
    ,(3535,1)
    ,(3554,0)
    ,(3567,1)
    ,(3599,1)
    ,(3603,0)
    ,(3605,0)
    ,(3617,0)
    ,(3654,0)
    ,(3661,0)
    ,(3666,1)
    ,(3696,1)
    ,(3703,1)
    ,(3707,1)
    ,(3719,1)
    ,(3723,0)
    ,(3728,1)
    ,(3768,1)
    ,(3784,0)
    ,(3796,1)
    ,(3821,1)
    ,(3829,1)
    ,(3852,0)
    ,(3878,1)
    ,(3893,0)
    ,(3900,0)
    ,(3904,1)
    ,(3911,1)
    ,(3928,1)
    ,(3934,1)
    ,(3939,1)
    ,(3956,0)
    ,(3959,1)
    ,(3968,0)
    ,(3976,1)
    ,(3983,1)
    ,(3986,1)
    ,(4002,1)
    ,(4016,1)
    ,(4031,1)
    ,(4078,1)
    ,(4091,1)
    ,(4136,1)
    ,(4165,1)
    ,(4200,1)
    ,(4210,1)
    ,(4218,1)
    ,(4221,1)
    ,(4244,1)
    ,(4276,1)
    ,(4291,1)
    ,(4309,0)
    ,(4313,1)
    ,(4323,1)
    ,(4324,1)
    ,(4330,1)
    ,(4348,1)
    ,(4352,1)
    ,(4356,1)
    ,(4370,1)
    ,(4371,1)
    ,(4382,0)
    ,(4385,1)
    ,(4417,1)
    ,(4440,1)
    ,(4442,1)
    ,(4450,1)
    ,(4476,0)
    ,(4488,1)
    ,(4519,0)
    ,(4533,0)
    ,(4550,1)
    ,(4554,1)
    ,(4560,1)
    ,(4578,0)
    ,(4584,0)
    ,(4628,1)
    ,(4631,1)
    ,(4637,1)
    ,(4639,1)
    ,(4644,0)
    ,(4666,1)
    ,(4722,1)
    ,(4736,0)
    ,(4746,1)
    ,(4750,1)
    ,(4752,0)
    ,(4774,1)
    ,(4840,1)
    ,(4856,1)
    ,(4885,0)
    ,(4910,0)
    ,(4950,1)
    ,(4952,0)
    ,(5000,0)
    ,(5008,1)
    ,(5031,1)
    ,(5033,0)
    ,(5060,0)
    ,(5074,1)
    ,(5080,1)
    ,(5105,1)
    ,(5159,1)
    ,(5179,0)
    ,(5199,1)
    ,(5215,1)
    ,(5225,1)
    ,(5231,1)
    ,(5252,0)
    ,(5262,1)
    ,(5270,1)
    ,(5275,0)
    ,(5293,1)
    ,(5337,1)
    ,(5377,1)
    ,(5407,1)
    ,(5410,1)
    ,(5444,0)
    ,(5450,1)
    ,(5484,1)
    ,(5490,0)
    ,(5491,1)
    ,(5525,1)
    ,(5544,0)
    ,(5545,1)
    ,(5550,1)
    ,(5559,1)
    ,(5560,0)
    ,(5586,1)
    ,(5644,1)
    ,(5661,1)
    ,(5673,0)
    ,(5676,1)
    ,(5679,1)
    ,(5735,0)
    ,(5798,1)
    ,(5800,0)
    ,(5853,1)
    ,(5865,0)
    ,(5868,0)
    ,(5902,1)
    ,(5919,0)
    ,(5925,1)
    ,(5961,1)
    ,(5968,1)
    ,(5977,0)
    ,(5987,0)
    ,(6054,1)
    ,(6059,1)
    ,(6103,1)
    ,(6104,0)
    ,(6125,1)
    ,(6162,0)
    ,(6180,1)
    ,(6184,0)
    ,(6204,1)
    ,(6219,1)
    ,(6225,0)
    ,(6279,1)
    ,(6314,1)
    ,(6322,1)
    ,(6331,0)
    ,(6391,0)
    ,(6402,0)
    ,(6407,1)
    ,(6414,1)
    ,(6436,0)
    ,(6493,1)
    ,(6513,0)
    ,(6517,1)
    ,(6522,0)
    ,(6531,0)
    ,(6533,1)
    ,(6543,1)
    ,(6568,0)
    ,(6569,1)
    ,(6597,1)
    ,(6621,1)
    ,(6639,1)
    ,(6678,1)
    ,(6693,1)
    ,(6703,0)
    ,(6707,1)
    ,(6709,1)
    ,(6713,1)
    ,(6723,1)
    ,(6756,1)
    ,(6768,1)
    ,(6778,0)
    ,(6801,0)
    ,(6809,0)
    ,(6810,1)
    ,(6811,1)
    ,(6874,1)
    ,(6883,0)
    ,(6886,0)
    ,(6897,1)
    ,(6920,1)
    ,(6998,1)
    ,(7004,0)
    ,(7033,1)
    ,(7063,0)
    ,(7068,1)
    ,(7070,1)
    ,(7074,0)
    ,(7084,1)
    ,(7090,1)
    ,(7095,1)
    ,(7158,1)
    ,(7168,1)
    ,(7175,0)
    ,(7179,0)
    ,(7198,0)
    ,(7203,0)
    ,(7251,0)
    ,(7316,1)
    ,(7321,1)
    ,(7346,1)
    ,(7372,1)
    ,(7374,1)
    ,(7415,1)
    ,(7467,0)
    ,(7487,1)
    ,(7520,0)
    ,(7529,1)
    ,(7544,0)
    ,(7553,0)
    ,(7573,1)
    ,(7600,1)
    ,(7628,1)
    ,(7653,1)
    ,(7655,0)
    ,(7671,1)
    ,(7673,1)
    ,(7711,0)
    ,(7717,1)
)

foreach ($pair in $changes) {
    $row = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($row, 2).Value = $val
}
